$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 515.75
$ws.Cells.Item(32, 9).Value = 552.2
$ws.Cells.Item(32, 10).Value = 455
$ws.Cells.Item(32, 11).Value = 552.2
$ws.Cells.Item(32, 12).Value = 455
$ws.Cells.Item(32, 13).Value = -226.2
$ws.Cells.Item(32, 14).Value = -1107
$ws.Cells.Item(137, 8).Value = 1952.1951
$ws.Cells.Item(137, 9).Value = 1035.4
$ws.Cells.Item(137, 10).Value = 2825.3333
$ws.Cells.Item(137, 11).Value = 3106.2
$ws.Cells.Item(137, 12).Value = 8475.999899999999
$ws.Cells.Item(137, 13).Value = -556.2000000000003
$ws.Cells.Item(137, 14).Value = -13575.9999
$ws.Cells.Item(138, 8).Value = 1679.3636
$ws.Cells.Item(138, 9).Value = 898.6326
$ws.Cells.Item(138, 10).Value = 2660.282
$ws.Cells.Item(138, 11).Value = 2695.8978
$ws.Cells.Item(138, 12).Value = 7980.846
$ws.Cells.Item(138, 13).Value = 2444.1022
$ws.Cells.Item(138, 14).Value = -18260.846
$ws.Cells.Item(141, 8).Value = 736.14636
$ws.Cells.Item(141, 9).Value = 636.5526
$ws.Cells.Item(141, 10).Value = 1997.6666
$ws.Cells.Item(141, 11).Value = 1909.6578
$ws.Cells.Item(141, 12).Value = 5992.9998
$ws.Cells.Item(141, 13).Value = 3270.3422
$ws.Cells.Item(141, 14).Value = -16352.9998
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8198922
$ws.Cells.Item(32, 9).Value = 9092741
$ws.Cells.Item(32, 10).Value = 5578.8335
$ws.Cells.Item(32, 11).Value = 9092741
$ws.Cells.Item(32, 12).Value = 5578.8335
$ws.Cells.Item(32, 13).Value = -9092454
$ws.Cells.Item(32, 14).Value = -6152.8335
$ws.Cells.Item(74, 8).Value = 4810187.5
$ws.Cells.Item(74, 9).Value = 5683407
$ws.Cells.Item(74, 10).Value = 7479.75
$ws.Cells.Item(74, 11).Value = 5683407
$ws.Cells.Item(74, 12).Value = 7479.75
$ws.Cells.Item(74, 13).Value = -5682533
$ws.Cells.Item(74, 14).Value = -9227.75
$ws.Cells.Item(77, 8).Value = 4810187.5
$ws.Cells.Item(77, 9).Value = 5683407
$ws.Cells.Item(77, 10).Value = 7479.75
$ws.Cells.Item(77, 11).Value = 28417035
$ws.Cells.Item(77, 12).Value = 37398.75
$ws.Cells.Item(77, 13).Value = -28412667
$ws.Cells.Item(77, 14).Value = -46134.75
$ws.Cells.Item(97, 8).Value = 978.2308
$ws.Cells.Item(97, 9).Value = 1245.5714
$ws.Cells.Item(97, 10).Value = 666.3333
$ws.Cells.Item(97, 11).Value = 1245.5714
$ws.Cells.Item(97, 12).Value = 666.3333
$ws.Cells.Item(97, 13).Value = -749.5714
$ws.Cells.Item(97, 14).Value = -1658.3333
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2514.8696
$ws.Cells.Item(3, 9).Value = 2985.647
$ws.Cells.Item(3, 10).Value = 1181
$ws.Cells.Item(3, 11).Value = 2985.647
$ws.Cells.Item(3, 12).Value = 1181
$ws.Cells.Item(3, 13).Value = -2871.647
$ws.Cells.Item(3, 14).Value = -1409
$ws.Cells.Item(80, 8).Value = 7685.8335
$ws.Cells.Item(80, 9).Value = 11890
$ws.Cells.Item(80, 10).Value = 4682.857
$ws.Cells.Item(80, 11).Value = 11890
$ws.Cells.Item(80, 12).Value = 4682.857
$ws.Cells.Item(80, 13).Value = -10892
$ws.Cells.Item(80, 14).Value = -6678.857
$ws.Cells.Item(83, 8).Value = 7685.8335
$ws.Cells.Item(83, 9).Value = 11890
$ws.Cells.Item(83, 10).Value = 4682.857
$ws.Cells.Item(83, 11).Value = 59450
$ws.Cells.Item(83, 12).Value = 23414.285
$ws.Cells.Item(83, 13).Value = -54458
$ws.Cells.Item(83, 14).Value = -33398.285
$ws.Cells.Item(94, 8).Value = 1098.5862
$ws.Cells.Item(94, 9).Value = 880
$ws.Cells.Item(94, 10).Value = 2993
$ws.Cells.Item(94, 11).Value = 880
$ws.Cells.Item(94, 12).Value = 2993
$ws.Cells.Item(94, 13).Value = -429
$ws.Cells.Item(94, 14).Value = -3895
$ws.Cells.Item(105, 8).Value = 2710.7778
$ws.Cells.Item(105, 9).Value = 2710.7778
$ws.Cells.Item(105, 11).Value = 2710.7778
$ws.Cells.Item(105, 13).Value = -963.7777999999998
$ws.Cells.Item(107, 8).Value = 3107.8235
$ws.Cells.Item(107, 9).Value = 2949.5454
$ws.Cells.Item(107, 11).Value = 2949.5454
$ws.Cells.Item(107, 13).Value = -1029.5454
$ws.Cells.Item(134, 8).Value = 771762.75
$ws.Cells.Item(134, 9).Value = 1135493.2
$ws.Cells.Item(134, 11).Value = 3406479.6
$ws.Cells.Item(134, 13).Value = -3403944.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 104595.945
$ws.Cells.Item(31, 9).Value = 199112.62
$ws.Cells.Item(31, 10).Value = 25002.947
$ws.Cells.Item(31, 11).Value = 199112.62
$ws.Cells.Item(31, 12).Value = 25002.947
$ws.Cells.Item(31, 13).Value = -198817.62
$ws.Cells.Item(31, 14).Value = -25592.947
$ws.Cells.Item(34, 8).Value = 104595.945
$ws.Cells.Item(34, 9).Value = 199112.62
$ws.Cells.Item(34, 10).Value = 25002.947
$ws.Cells.Item(34, 11).Value = 199112.62
$ws.Cells.Item(34, 12).Value = 25002.947
$ws.Cells.Item(34, 13).Value = -198910.62
$ws.Cells.Item(34, 14).Value = -25406.947
$ws.Cells.Item(58, 8).Value = 175968.4
$ws.Cells.Item(58, 9).Value = 269915.8
$ws.Cells.Item(58, 10).Value = 3105.16
$ws.Cells.Item(58, 11).Value = 269915.8
$ws.Cells.Item(58, 12).Value = 3105.16
$ws.Cells.Item(58, 13).Value = -269712.8
$ws.Cells.Item(58, 14).Value = -3511.16
$ws.Cells.Item(86, 8).Value = 3330.5557
$ws.Cells.Item(86, 9).Value = 2961.2
$ws.Cells.Item(86, 10).Value = 3792.25
$ws.Cells.Item(86, 11).Value = 2961.2
$ws.Cells.Item(86, 12).Value = 3792.25
$ws.Cells.Item(86, 13).Value = -1838.2
$ws.Cells.Item(86, 14).Value = -6038.25
$ws.Cells.Item(89, 8).Value = 3330.5557
$ws.Cells.Item(89, 9).Value = 2961.2
$ws.Cells.Item(89, 10).Value = 3792.25
$ws.Cells.Item(89, 11).Value = 14806
$ws.Cells.Item(89, 12).Value = 18961.25
$ws.Cells.Item(89, 13).Value = -9190
$ws.Cells.Item(89, 14).Value = -30193.25
$ws.Cells.Item(132, 8).Value = 5564739.5
$ws.Cells.Item(132, 9).Value = 11592.728
$ws.Cells.Item(132, 10).Value = 20835894
$ws.Cells.Item(132, 11).Value = 34778.18399999999
$ws.Cells.Item(132, 12).Value = 62507682
$ws.Cells.Item(132, 13).Value = -32248.18399999999
$ws.Cells.Item(132, 14).Value = -62512742
$ws.Cells.Item(134, 8).Value = 6301.0293
$ws.Cells.Item(134, 9).Value = 8397.234
$ws.Cells.Item(134, 11).Value = 25191.702
$ws.Cells.Item(134, 13).Value = -22656.702
$ws.Cells.Item(136, 8).Value = 175968.4
$ws.Cells.Item(136, 9).Value = 269915.8
$ws.Cells.Item(136, 10).Value = 3105.16
$ws.Cells.Item(136, 11).Value = 809747.3999999999
$ws.Cells.Item(136, 12).Value = 9315.48
$ws.Cells.Item(136, 13).Value = -807197.3999999999
$ws.Cells.Item(136, 14).Value = -14415.48
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 2928717.2
$ws.Cells.Item(7, 9).Value = 2000176.8
$ws.Cells.Item(7, 10).Value = 5250069
$ws.Cells.Item(7, 11).Value = 6000530.4
$ws.Cells.Item(7, 12).Value = 15750207
$ws.Cells.Item(7, 13).Value = -6000418.4
$ws.Cells.Item(7, 14).Value = -15750431
$ws.Cells.Item(86, 8).Value = 758.75
$ws.Cells.Item(86, 9).Value = 820
$ws.Cells.Item(86, 10).Value = 656.6667
$ws.Cells.Item(86, 11).Value = 2460
$ws.Cells.Item(86, 12).Value = 1970.0001
$ws.Cells.Item(86, 13).Value = -1274
$ws.Cells.Item(86, 14).Value = -4342.0001
$ws.Cells.Item(89, 8).Value = 758.75
$ws.Cells.Item(89, 9).Value = 820
$ws.Cells.Item(89, 10).Value = 656.6667
$ws.Cells.Item(89, 11).Value = 7380
$ws.Cells.Item(89, 12).Value = 5910.0003
$ws.Cells.Item(89, 13).Value = -1452
$ws.Cells.Item(89, 14).Value = -17766.0003
$ws.Cells.Item(117, 8).Value = 2591.4707
$ws.Cells.Item(117, 9).Value = 1243.1333
$ws.Cells.Item(117, 10).Value = 3655.9473
$ws.Cells.Item(117, 11).Value = 3729.3999
$ws.Cells.Item(117, 12).Value = 10967.8419
$ws.Cells.Item(117, 13).Value = -287.3998999999999
$ws.Cells.Item(117, 14).Value = -17851.8419
$ws.Cells.Item(131, 8).Value = 6567.9316
$ws.Cells.Item(131, 9).Value = 643.36365
$ws.Cells.Item(131, 10).Value = 8542.788
$ws.Cells.Item(131, 11).Value = 1930.09095
$ws.Cells.Item(131, 12).Value = 25628.364
$ws.Cells.Item(131, 13).Value = 3109.90905
$ws.Cells.Item(131, 14).Value = -35708.364
$ws.Cells.Item(137, 8).Value = 3866.75
$ws.Cells.Item(137, 10).Value = 5619.0835
$ws.Cells.Item(137, 12).Value = 16857.2505
$ws.Cells.Item(137, 14).Value = -27057.2505
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 137.90909
$ws.Cells.Item(2, 9).Value = 48.833332
$ws.Cells.Item(2, 10).Value = 244.8
$ws.Cells.Item(2, 11).Value = 48.833332
$ws.Cells.Item(2, 12).Value = 244.8
$ws.Cells.Item(2, 13).Value = 64.166668
$ws.Cells.Item(2, 14).Value = -470.8
$ws.Cells.Item(48, 8).Value = 23000
$ws.Cells.Item(48, 9).Value = 23000
$ws.Cells.Item(48, 11).Value = 23000
$ws.Cells.Item(48, 13).Value = -22515
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3026.5806
$ws.Cells.Item(40, 9).Value = 2914.8096
$ws.Cells.Item(40, 10).Value = 3261.3
$ws.Cells.Item(40, 11).Value = 2914.8096
$ws.Cells.Item(40, 12).Value = 3261.3
$ws.Cells.Item(40, 13).Value = -2778.8096
$ws.Cells.Item(40, 14).Value = -3533.3
$ws.Cells.Item(132, 8).Value = 990119.5
$ws.Cells.Item(132, 9).Value = 1117618.2
$ws.Cells.Item(132, 10).Value = 2004.75
$ws.Cells.Item(132, 11).Value = 3352854.6
$ws.Cells.Item(132, 12).Value = 6014.25
$ws.Cells.Item(132, 13).Value = -3350324.6
$ws.Cells.Item(132, 14).Value = -11074.25
$ws.Cells.Item(136, 8).Value = 27883.174
$ws.Cells.Item(136, 9).Value = 3538.1592
$ws.Cells.Item(136, 10).Value = 161780.75
$ws.Cells.Item(136, 11).Value = 10614.4776
$ws.Cells.Item(136, 12).Value = 485342.25
$ws.Cells.Item(136, 13).Value = -8064.4776
$ws.Cells.Item(136, 14).Value = -490442.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1878.6
$ws.Cells.Item(81, 9).Value = 1999.6666
$ws.Cells.Item(81, 10).Value = 1697
$ws.Cells.Item(81, 11).Value = 3999.3332
$ws.Cells.Item(81, 12).Value = 3394
$ws.Cells.Item(81, 13).Value = -2938.3332
$ws.Cells.Item(81, 14).Value = -5516
$ws.Cells.Item(84, 8).Value = 1878.6
$ws.Cells.Item(84, 9).Value = 1999.6666
$ws.Cells.Item(84, 10).Value = 1697
$ws.Cells.Item(84, 11).Value = 19996.666
$ws.Cells.Item(84, 12).Value = 16970
$ws.Cells.Item(84, 13).Value = -14692.666
$ws.Cells.Item(84, 14).Value = -27578
$ws.Cells.Item(122, 8).Value = 1747.2642
$ws.Cells.Item(122, 9).Value = 1682.4584
$ws.Cells.Item(122, 10).Value = 2369.4
$ws.Cells.Item(122, 11).Value = 5047.3752
$ws.Cells.Item(122, 12).Value = 7108.200000000001
$ws.Cells.Item(122, 13).Value = -2597.3752
$ws.Cells.Item(122, 14).Value = -12008.2
$ws.Cells.Item(126, 8).Value = 2870
$ws.Cells.Item(126, 9).Value = 2857
$ws.Cells.Item(126, 11).Value = 8571
$ws.Cells.Item(126, 13).Value = -6101
$ws.Cells.Item(136, 8).Value = 7762600.5
$ws.Cells.Item(136, 9).Value = 10269944
$ws.Cells.Item(136, 10).Value = 31624.166
$ws.Cells.Item(136, 11).Value = 30809832
$ws.Cells.Item(136, 12).Value = 94872.49800000001
$ws.Cells.Item(136, 13).Value = -30807282
$ws.Cells.Item(136, 14).Value = -99972.49800000001
